# The source data rows were re-sorted: the observation records that used to
# sit in rows 75-78 have been rotated among those same four rows (the row
# index itself doesn't move, only the record contents do).
#   new row 75 <- old row 77
#   new row 76 <- old row 75
#   new row 77 <- old row 78
#   new row 78 <- old row 76
# Only the columns whose content actually differs between these four
# records need to be touched: A,B,D,E,F,G,H,I,J,K,L,Q,R,AC.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","D","E","F","G","H","I","J","K","L","Q","R","AC")
$rows = @(75, 76, 77, 78)

# Snapshot the cells we are about to rewrite (read Value2 - the raw
# number/string/bool, no formatting applied).
$snapshot = @{}
foreach ($r in $rows) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Destination row -> source row (cyclic rotation of the 4 records).
$sourceFor = @{ 75 = 77; 76 = 75; 77 = 78; 78 = 76 }

foreach ($destRow in $rows) {
    $srcRow = $sourceFor[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $val = $srcVals[$col]
        $cell = $ws.Range("$col$destRow")
        if ($null -eq $val) {
            $cell.Value = $null
        } elseif ($col -eq "I") {
            # "Antal" (count) is stored as text ("1"), not a number - force
            # text so Excel doesn't silently coerce it to a numeric cell.
            $cell.NumberFormat = "@"
            $cell.Value = $val
        } else {
            $cell.Value = $val
        }
    }
}
